# Update the "Datos" worksheet with new test credentials and refresh the
# active cell selection, matching the upstream project restructuring commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# nit (G2): 123456789 -> 800270811 (leading apostrophe keeps it stored as
# text with the existing quote-prefixed text style, avoiding a reformat)
$ws.Range("G2").Value = "'800270811"

# username (H2): usuariotest -> autouser1994
$ws.Range("H2").Value = "'autouser1994"

# Update the saved selection/active cell on the sheet to I3.
$ws.Activate()
$ws.Range("I3").Select()
